$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 69.319552
$ws.Range("H2").Value = 207.958656
$ws.Range("I2").Value = 0.6721843675940576
$ws.Range("J2").Value = 0.6721843675940575
$ws.Range("M2").Value = 119.0164006666667
$ws.Range("N2").Value = 357.049202
$ws.Range("O2").Value = 0.9176278005170622
$ws.Range("P2").Value = 0.9176278005170622
$ws.Range("Q2").Value = 8250.163574865834
$ws.Range("R2").Value = 74251.47217379251
$ws.Range("S2").Value = 0.6168150627772875
$ws.Range("T2").Value = 0.6168150627772874

$ws.Range("G3").Value = 69.319552
$ws.Range("H3").Value = 207.958656
$ws.Range("I3").Value = 0.6721843675940576
$ws.Range("J3").Value = 0.6721843675940575
$ws.Range("O3").Value = 0.001755838010330732
$ws.Range("P3").Value = 0.001755838010330731
$ws.Range("Q3").Value = 15.78630332258133
$ws.Range("R3").Value = 142.076729903232
$ws.Range("S3").Value = 0.001180246862571771
$ws.Range("T3").Value = 0.001180246862571771

$ws.Range("G4").Value = 69.319552
$ws.Range("H4").Value = 207.958656
$ws.Range("I4").Value = 0.6721843675940576
$ws.Range("J4").Value = 0.6721843675940575
$ws.Range("M4").Value = 7.816301333333333
$ws.Range("N4").Value = 23.448904
$ws.Range("O4").Value = 0.06026442877207647
$ws.Range("P4").Value = 0.06026442877207646
$ws.Range("Q4").Value = 541.8225067236693
$ws.Range("R4").Value = 4876.402560513025
$ws.Range("S4").Value = 0.04050880694257535
$ws.Range("T4").Value = 0.04050880694257534

$ws.Range("G5").Value = 69.319552
$ws.Range("H5").Value = 207.958656
$ws.Range("I5").Value = 0.6721843675940576
$ws.Range("J5").Value = 0.6721843675940575
$ws.Range("M5").Value = 0.105045
$ws.Range("N5").Value = 0.315135
$ws.Range("O5").Value = 0.000809906968832672
$ws.Range("P5").Value = 0.000809906968832672
$ws.Range("Q5").Value = 7.28167233984
$ws.Range("R5").Value = 65.53505105856
$ws.Range("S5").Value = 0.0005444068036548098
$ws.Range("T5").Value = 0.0005444068036548097

$ws.Range("G6").Value = 69.319552
$ws.Range("H6").Value = 207.958656
$ws.Range("I6").Value = 0.6721843675940576
$ws.Range("J6").Value = 0.6721843675940575
$ws.Range("M6").Value = 2.534602333333333
$ws.Range("N6").Value = 7.603807
$ws.Range("O6").Value = 0.01954202573169801
$ws.Range("P6").Value = 0.01954202573169801
$ws.Range("Q6").Value = 175.6974982448213
$ws.Range("R6").Value = 1581.277484203392
$ws.Range("S6").Value = 0.01313584420796823
$ws.Range("T6").Value = 0.01313584420796823

$ws.Range("I7").Value = 0.1094364499261462
$ws.Range("J7").Value = 0.1094364499261462
$ws.Range("M7").Value = 119.0164006666667
$ws.Range("N7").Value = 357.049202
$ws.Range("O7").Value = 0.9176278005170622
$ws.Range("P7").Value = 0.9176278005170622
$ws.Range("Q7").Value = 1343.185971692481
$ws.Range("R7").Value = 12088.67374523233
$ws.Range("S7").Value = 0.1004219288421251
$ws.Range("T7").Value = 0.1004219288421251

$ws.Range("I8").Value = 0.1094364499261462
$ws.Range("J8").Value = 0.1094364499261462
$ws.Range("O8").Value = 0.001755838010330732
$ws.Range("P8").Value = 0.001755838010330731
$ws.Range("S8").Value = 0.0001921526784959833
$ws.Range("T8").Value = 0.0001921526784959832

$ws.Range("I9").Value = 0.1094364499261462
$ws.Range("J9").Value = 0.1094364499261462
$ws.Range("M9").Value = 7.816301333333333
$ws.Range("N9").Value = 23.448904
$ws.Range("O9").Value = 0.06026442877207647
$ws.Range("P9").Value = 0.06026442877207646
$ws.Range("Q9").Value = 88.21260131079553
$ws.Range("R9").Value = 793.9134117971598
$ws.Range("S9").Value = 0.00659512514164315
$ws.Range("T9").Value = 0.006595125141643147

$ws.Range("I10").Value = 0.1094364499261462
$ws.Range("J10").Value = 0.1094364499261462
$ws.Range("M10").Value = 0.105045
$ws.Range("N10").Value = 0.315135
$ws.Range("O10").Value = 0.000809906968832672
$ws.Range("P10").Value = 0.000809906968832672
$ws.Range("Q10").Value = 1.185508632475
$ws.Range("R10").Value = 10.669577692275
$ws.Range("S10").Value = 0.00008863334343949355
$ws.Range("T10").Value = 0.00008863334343949352

$ws.Range("I11").Value = 0.1094364499261462
$ws.Range("J11").Value = 0.1094364499261462
$ws.Range("M11").Value = 2.534602333333333
$ws.Range("N11").Value = 7.603807
$ws.Range("O11").Value = 0.01954202573169801
$ws.Range("P11").Value = 0.01954202573169801
$ws.Range("Q11").Value = 28.60481646968389
$ws.Range("R11").Value = 257.443348227155
$ws.Range("S11").Value = 0.00213860992044243
$ws.Range("T11").Value = 0.002138609920442429

$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.03495566666666667
$ws.Range("H12").Value = 0.104867
$ws.Range("I12").Value = 0.0003389614043114707
$ws.Range("J12").Value = 0.0003389614043114706
$ws.Range("M12").Value = 119.0164006666667
$ws.Range("N12").Value = 357.049202
$ws.Range("O12").Value = 0.9176278005170622
$ws.Range("P12").Value = 0.9176278005170622
$ws.Range("Q12").Value = 4.160297629570445
$ws.Range("R12").Value = 37.442678666134
$ws.Range("S12").Value = 0.0003110404078985095
$ws.Range("T12").Value = 0.0003110404078985094

$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.03495566666666667
$ws.Range("H13").Value = 0.104867
$ws.Range("I13").Value = 0.0003389614043114707
$ws.Range("J13").Value = 0.0003389614043114706
$ws.Range("O13").Value = 0.001755838010330732
$ws.Range("P13").Value = 0.001755838010330731
$ws.Range("Q13").Value = 0.007960535533222223
$ws.Range("R13").Value = 0.071644819799
$ws.Range("S13").Value = 0.0000005951613177251633
$ws.Range("T13").Value = 0.000000595161317725163

$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.03495566666666667
$ws.Range("H14").Value = 0.104867
$ws.Range("I14").Value = 0.0003389614043114707
$ws.Range("J14").Value = 0.0003389614043114706
$ws.Range("M14").Value = 7.816301333333333
$ws.Range("N14").Value = 23.448904
$ws.Range("O14").Value = 0.06026442877207647
$ws.Range("P14").Value = 0.06026442877207646
$ws.Range("Q14").Value = 0.2732240239742222
$ws.Range("R14").Value = 2.459016215768
$ws.Range("S14").Value = 0.00002042731540661164
$ws.Range("T14").Value = 0.00002042731540661163

$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.03495566666666667
$ws.Range("H15").Value = 0.104867
$ws.Range("I15").Value = 0.0003389614043114707
$ws.Range("J15").Value = 0.0003389614043114706
$ws.Range("M15").Value = 0.105045
$ws.Range("N15").Value = 0.315135
$ws.Range("O15").Value = 0.000809906968832672
$ws.Range("P15").Value = 0.000809906968832672
$ws.Range("Q15").Value = 0.003671918005
$ws.Range("R15").Value = 0.033047262045
$ws.Range("S15").Value = 0.000000274527203517169
$ws.Range("T15").Value = 0.0000002745272035171689

$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.03495566666666667
$ws.Range("H16").Value = 0.104867
$ws.Range("I16").Value = 0.0003389614043114707
$ws.Range("J16").Value = 0.0003389614043114706
$ws.Range("M16").Value = 2.534602333333333
$ws.Range("N16").Value = 7.603807
$ws.Range("O16").Value = 0.01954202573169801
$ws.Range("P16").Value = 0.01954202573169801
$ws.Range("Q16").Value = 0.08859871429655557
$ws.Range("R16").Value = 0.797388428669
$ws.Range("S16").Value = 0.000006623992485107253
$ws.Range("T16").Value = 0.000006623992485107251

$ws.Range("G17").Value = 22.33376
$ws.Range("H17").Value = 67.00128000000001
$ws.Range("I17").Value = 0.2165681096957675
$ws.Range("J17").Value = 0.2165681096957675
$ws.Range("M17").Value = 119.0164006666667
$ws.Range("N17").Value = 357.049202
$ws.Range("O17").Value = 0.9176278005170622
$ws.Range("P17").Value = 0.9176278005170622
$ws.Range("Q17").Value = 2658.083728553173
$ws.Range("R17").Value = 23922.75355697856
$ws.Range("S17").Value = 0.198728918162265
$ws.Range("T17").Value = 0.198728918162265

$ws.Range("G18").Value = 22.33376
$ws.Range("H18").Value = 67.00128000000001
$ws.Range("I18").Value = 0.2165681096957675
$ws.Range("J18").Value = 0.2165681096957675
$ws.Range("O18").Value = 0.001755838010330732
$ws.Range("P18").Value = 0.001755838010330731
$ws.Range("Q18").Value = 5.086119276906667
$ws.Range("R18").Value = 45.77507349216001
$ws.Range("S18").Value = 0.0003802585188293041
$ws.Range("T18").Value = 0.0003802585188293039

$ws.Range("G19").Value = 22.33376
$ws.Range("H19").Value = 67.00128000000001
$ws.Range("I19").Value = 0.2165681096957675
$ws.Range("J19").Value = 0.2165681096957675
$ws.Range("M19").Value = 7.816301333333333
$ws.Range("N19").Value = 23.448904
$ws.Range("O19").Value = 0.06026442877207647
$ws.Range("P19").Value = 0.06026442877207646
$ws.Range("Q19").Value = 174.5673980663467
$ws.Range("R19").Value = 1571.10658259712
$ws.Range("S19").Value = 0.01305135342106383
$ws.Range("T19").Value = 0.01305135342106382

$ws.Range("G20").Value = 22.33376
$ws.Range("H20").Value = 67.00128000000001
$ws.Range("I20").Value = 0.2165681096957675
$ws.Range("J20").Value = 0.2165681096957675
$ws.Range("M20").Value = 0.105045
$ws.Range("N20").Value = 0.315135
$ws.Range("O20").Value = 0.000809906968832672
$ws.Range("P20").Value = 0.000809906968832672
$ws.Range("Q20").Value = 2.3460498192
$ws.Range("R20").Value = 21.1144483728
$ws.Range("S20").Value = 0.0001754000212695207
$ws.Range("T20").Value = 0.0001754000212695206

$ws.Range("G21").Value = 22.33376
$ws.Range("H21").Value = 67.00128000000001
$ws.Range("I21").Value = 0.2165681096957675
$ws.Range("J21").Value = 0.2165681096957675
$ws.Range("M21").Value = 2.534602333333333
$ws.Range("N21").Value = 7.603807
$ws.Range("O21").Value = 0.01954202573169801
$ws.Range("P21").Value = 0.01954202573169801
$ws.Range("Q21").Value = 56.60720020810668
$ws.Range("R21").Value = 509.46480187296
$ws.Range("S21").Value = 0.004232179572339887
$ws.Range("T21").Value = 0.004232179572339886

$ws.Range("G22").Value = 0.1518126666666667
$ws.Range("H22").Value = 0.455438
$ws.Range("I22").Value = 0.001472111379717238
$ws.Range("J22").Value = 0.001472111379717237
$ws.Range("M22").Value = 119.0164006666667
$ws.Range("N22").Value = 357.049202
$ws.Range("O22").Value = 0.9176278005170622
$ws.Range("P22").Value = 0.9176278005170622
$ws.Range("Q22").Value = 18.06819716227511
$ws.Range("R22").Value = 162.613774460476
$ws.Range("S22").Value = 0.001350850327486067
$ws.Range("T22").Value = 0.001350850327486066

$ws.Range("G23").Value = 0.1518126666666667
$ws.Range("H23").Value = 0.455438
$ws.Range("I23").Value = 0.001472111379717238
$ws.Range("J23").Value = 0.001472111379717237
$ws.Range("O23").Value = 0.001755838010330732
$ws.Range("P23").Value = 0.001755838010330731
$ws.Range("Q23").Value = 0.03457265280955556
$ws.Range("R23").Value = 0.311153875286
$ws.Range("S23").Value = 0.000002584789115947943
$ws.Range("T23").Value = 0.000002584789115947942

$ws.Range("G24").Value = 0.1518126666666667
$ws.Range("H24").Value = 0.455438
$ws.Range("I24").Value = 0.001472111379717238
$ws.Range("J24").Value = 0.001472111379717237
$ws.Range("M24").Value = 7.816301333333333
$ws.Range("N24").Value = 23.448904
$ws.Range("O24").Value = 0.06026442877207647
$ws.Range("P24").Value = 0.06026442877207646
$ws.Range("Q24").Value = 1.186613548883556
$ws.Range("R24").Value = 10.679521939952
$ws.Range("S24").Value = 0.00008871595138753269
$ws.Range("T24").Value = 0.00008871595138753266

$ws.Range("G25").Value = 0.1518126666666667
$ws.Range("H25").Value = 0.455438
$ws.Range("I25").Value = 0.001472111379717238
$ws.Range("J25").Value = 0.001472111379717237
$ws.Range("M25").Value = 0.105045
$ws.Range("N25").Value = 0.315135
$ws.Range("O25").Value = 0.000809906968832672
$ws.Range("P25").Value = 0.000809906968832672
$ws.Range("Q25").Value = 0.01594716157
$ws.Range("R25").Value = 0.14352445413
$ws.Range("S25").Value = 0.000001192273265330871
$ws.Range("T25").Value = 0.00000119227326533087

$ws.Range("G26").Value = 0.1518126666666667
$ws.Range("H26").Value = 0.455438
$ws.Range("I26").Value = 0.001472111379717238
$ws.Range("J26").Value = 0.001472111379717237
$ws.Range("M26").Value = 2.534602333333333
$ws.Range("N26").Value = 7.603807
$ws.Range("O26").Value = 0.01954202573169801
$ws.Range("P26").Value = 0.01954202573169801
$ws.Range("Q26").Value = 0.3847847391628889
$ws.Range("R26").Value = 3.463062652466
$ws.Range("S26").Value = 0.00002876803846235972
$ws.Range("T26").Value = 0.00002876803846235972
